{"js": "// SKILLS section clean-up:\n//  - \"Technical Skills:\" bullet  -> \"System Administration\" becomes\n//    \"Network Fundamentals\" (the only real wording change).\n//  - \"Project Management:\" and \"Digital Marketing & Content Creation:\"\n//    bullets keep their exact wording; their trailing sentence just gets\n//    consolidated back into a single run of text.\n// All three bullets are plain (non-bold) runs in Times New Roman, 12pt,\n// so a straightforward search-and-replace of the relevant text preserves\n// formatting while applying the edit.\n\nconst body = context.document.body;\n\n// 1) \"System Administration\" -> \"Network Fundamentals\" under Technical Skills.\nconst techHits = body.search(\"System Administration\", { matchCase: true });\ntechHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of techHits.items) {\n  hit.insertText(\"Network Fundamentals\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Re-unite the Project Management sentence into one contiguous run.\nconst pmHits = body.search(\n  \"Agile, Scrum, Problem Solving, Time Management, Negotiation, Organization, Flexibility\",\n  { matchCase: true }\n);\npmHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of pmHits.items) {\n  hit.insertText(\n    \"Agile, Scrum, Problem Solving, Time Management, Negotiation, Organization, Flexibility\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 3) Re-unite the Digital Marketing & Content Creation sentence into one run.\nconst dmHits = body.search(\n  \"Social Media Management, SEO, Content Strategy, Creativity\",\n  { matchCase: true }\n);\ndmHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of dmHits.items) {\n  hit.insertText(\n    \"Social Media Management, SEO, Content Strategy, Creativity\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n", "ps1": "# SKILLS section clean-up:\n#  - \"Technical Skills:\" bullet  -> \"System Administration\" becomes\n#    \"Network Fundamentals\" (the only real wording change).\n#  - \"Project Management:\" and \"Digital Marketing & Content Creation:\"\n#    bullets keep their exact wording; their trailing sentence just gets\n#    consolidated back into a single run of text.\n# All three bullets are plain (non-bold) runs in Times New Roman, 12pt,\n# so a straightforward Find/Replace over the relevant text preserves\n# formatting while applying the edit.\n\n$d = $word.ActiveDocument\n\n# wdFindContinue / wdReplaceAll style constants used below:\n#   Forward=$true, Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll)\n\n# 1) \"System Administration\" -> \"Network Fundamentals\" under Technical Skills.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"System Administration\", $false, $false, $false, $false, $false, $true, 1, $false, \"Network Fundamentals\", 2) | Out-Null\n\n# 2) Re-unite the Project Management sentence into one contiguous run.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"Agile, Scrum, Problem Solving, Time Management, Negotiation, Organization, Flexibility\", $false, $false, $false, $false, $false, $true, 1, $false, \"Agile, Scrum, Problem Solving, Time Management, Negotiation, Organization, Flexibility\", 2) | Out-Null\n\n# 3) Re-unite the Digital Marketing & Content Creation sentence into one run.\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Execute(\"Social Media Management, SEO, Content Strategy, Creativity\", $false, $false, $false, $false, $false, $true, 1, $false, \"Social Media Management, SEO, Content Strategy, Creativity\", 2) | Out-Null\n"}
